# "#12 Future works added on slides"
#
# Inserts a new "Future work" slide right after the "Conclusion" slide
# (position 7), using the same "Titre et contenu" (Title and Content)
# layout already used by the other content slides in this deck. All the
# slides that used to follow land one position further down, which is
# the only other observable change to the deck's slide order.

$p = $ppt.ActivePresentation

# "Titre et contenu" is CustomLayout #2 on this deck's master - the same
# layout used by the slides around it (Conclusion, annexe, ...).
$layouts = $p.SlideMaster.CustomLayouts
$newSlide = $p.Slides.AddSlide(7, $layouts.Item(2))

# Match the naming convention (French placeholder names) used throughout
# the rest of this presentation.
$titleShape = $newSlide.Shapes.Item(1)
$bodyShape = $newSlide.Shapes.Item(2)
$titleShape.Name = "Titre 1"
$bodyShape.Name = "Espace réservé du contenu 2"

$titleShape.TextFrame.TextRange.Text = "Future work"

$bodyText = "Virtual assistant launch every time the computer is up`r" + `
    "Generate a voice for the virtual assistant to communicate with the user`r" + `
    "Profiling each user that use the virtual assistant"
$bodyShape.TextFrame.TextRange.Text = $bodyText
